$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 4719.154  # H40: was 4732.0713
$ws.Cells.Item(40, 10).Value = 5674.875  # J40: was 5588.778
$ws.Cells.Item(40, 12).Value = 5674.875  # L40: was 5588.778
$ws.Cells.Item(40, 14).Value = -6024.875  # N40: was -5938.778
$ws.Cells.Item(48, 8).Value = 0  # H48: was 1550
$ws.Cells.Item(48, 9).Value = 0  # I48: was 1500
$ws.Cells.Item(48, 10).Value = 0  # J48: was 1600
$ws.Cells.Item(48, 11).Value = 0  # K48: was 4500
$ws.Cells.Item(48, 12).Value = 0  # L48: was 4800
$ws.Cells.Item(48, 13).Value = ""  # M48: was -4208
$ws.Cells.Item(48, 14).Value = 0  # N48: was -5384
$ws.Cells.Item(56, 8).Value = 0  # H56: was 1550
$ws.Cells.Item(56, 9).Value = 0  # I56: was 1500
$ws.Cells.Item(56, 10).Value = 0  # J56: was 1600
$ws.Cells.Item(56, 11).Value = 0  # K56: was 4500
$ws.Cells.Item(56, 12).Value = 0  # L56: was 4800
$ws.Cells.Item(56, 13).Value = ""  # M56: was -3966
$ws.Cells.Item(56, 14).Value = 0  # N56: was -5868
$ws.Cells.Item(64, 8).Value = 13498.5625  # H64: was 14284.071
$ws.Cells.Item(64, 9).Value = 9554.223  # I64: was 9998.286
$ws.Cells.Item(64, 11).Value = 9554.223  # K64: was 9998.286
$ws.Cells.Item(64, 13).Value = -9306.223  # M64: was -9750.286
$ws.Cells.Item(67, 8).Value = 13498.5625  # H67: was 14284.071
$ws.Cells.Item(67, 9).Value = 9554.223  # I67: was 9998.286
$ws.Cells.Item(67, 11).Value = 9554.223  # K67: was 9998.286
$ws.Cells.Item(67, 13).Value = -8696.223  # M67: was -9140.286
$ws.Cells.Item(125, 8).Value = 3089.5557  # H125: was 2878.5
$ws.Cells.Item(125, 9).Value = 1033  # I125: was 1300
$ws.Cells.Item(125, 10).Value = 4117.8335  # J125: was 3667.75
$ws.Cells.Item(125, 11).Value = 9297  # K125: was 11700
$ws.Cells.Item(125, 12).Value = 37060.5015  # L125: was 33009.75
$ws.Cells.Item(125, 13).Value = -6837  # M125: was -9240
$ws.Cells.Item(125, 14).Value = -41980.5015  # N125: was -37929.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(3, 8).Value = 2647.3333  # H3: was 2649
$ws.Cells.Item(3, 9).Value = 146  # I3: was 148.5
$ws.Cells.Item(3, 11).Value = 146  # K3: was 148.5
$ws.Cells.Item(3, 13).Value = -31  # M3: was -33.5
$ws.Cells.Item(8, 8).Value = 4012571  # H8: was 4012651
$ws.Cells.Item(8, 9).Value = 6683468.5  # I8: was 6683601.5
$ws.Cells.Item(8, 11).Value = 6683468.5  # K8: was 6683601.5
$ws.Cells.Item(8, 13).Value = -6683324.5  # M8: was -6683457.5
$ws.Cells.Item(61, 8).Value = 8249.75  # H61: was 6779.6
$ws.Cells.Item(61, 9).Value = 8666.333000000001  # I61: was 6724.5
$ws.Cells.Item(61, 11).Value = 8666.333000000001  # K61: was 6724.5
$ws.Cells.Item(61, 13).Value = -8454.333000000001  # M61: was -6512.5
$ws.Cells.Item(62, 8).Value = 0  # H62: was 5000
$ws.Cells.Item(62, 9).Value = 0  # I62: was 5000
$ws.Cells.Item(62, 11).Value = 0  # K62: was 5000
$ws.Cells.Item(62, 13).Value = ""  # M62: was -4376
$ws.Cells.Item(65, 8).Value = 0  # H65: was 5000
$ws.Cells.Item(65, 9).Value = 0  # I65: was 5000
$ws.Cells.Item(65, 11).Value = 0  # K65: was 15000
$ws.Cells.Item(65, 13).Value = ""  # M65: was -11880
$ws.Cells.Item(110, 8).Value = 1513.5555  # H110: was 1197.2307
$ws.Cells.Item(110, 9).Value = 1604.7142  # I110: was 1197.7273
$ws.Cells.Item(110, 11).Value = 1604.7142  # K110: was 1197.7273
$ws.Cells.Item(110, 13).Value = 440.2858000000001  # M110: was 847.2727
$ws.Cells.Item(132, 8).Value = 4089.4546  # H132: was 4178.8
$ws.Cells.Item(132, 9).Value = 3976  # I132: was 4073.5
$ws.Cells.Item(132, 11).Value = 11928  # K132: was 12220.5
$ws.Cells.Item(132, 13).Value = -9398  # M132: was -9690.5
$ws.Cells.Item(136, 8).Value = 8249.75  # H136: was 6779.6
$ws.Cells.Item(136, 9).Value = 8666.333000000001  # I136: was 6724.5
$ws.Cells.Item(136, 11).Value = 25998.999  # K136: was 20173.5
$ws.Cells.Item(136, 13).Value = -23448.999  # M136: was -17623.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(7, 8).Value = 3167103  # H7: was 3167853
$ws.Cells.Item(7, 10).Value = 500  # J7: was 5000
$ws.Cells.Item(7, 12).Value = 500  # L7: was 5000
$ws.Cells.Item(7, 14).Value = -726  # N7: was -5226

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 5437.8335  # H4: was 5508.75
$ws.Cells.Item(4, 9).Value = 1465.1428  # I4: was 1586.7142
$ws.Cells.Item(4, 11).Value = 1465.1428  # K4: was 1586.7142
$ws.Cells.Item(4, 13).Value = -1353.1428  # M4: was -1474.7142
$ws.Cells.Item(19, 8).Value = 972.7273  # H19: was 1068.5
$ws.Cells.Item(19, 9).Value = 70.09999999999999  # I19: was 76.22221999999999
$ws.Cells.Item(19, 11).Value = 70.09999999999999  # K19: was 76.22221999999999
$ws.Cells.Item(19, 13).Value = 99.90000000000001  # M19: was 93.77778000000001
$ws.Cells.Item(22, 8).Value = 1679.1875  # H22: was 962.44446
$ws.Cells.Item(22, 9).Value = 638  # I22: was 665.5
$ws.Cells.Item(22, 10).Value = 2026.25  # J22: was 1200
$ws.Cells.Item(22, 11).Value = 638  # K22: was 665.5
$ws.Cells.Item(22, 12).Value = 2026.25  # L22: was 1200
$ws.Cells.Item(22, 13).Value = -288  # M22: was -315.5
$ws.Cells.Item(22, 14).Value = -2726.25  # N22: was -1900
$ws.Cells.Item(24, 8).Value = 972.7273  # H24: was 1068.5
$ws.Cells.Item(24, 9).Value = 70.09999999999999  # I24: was 76.22221999999999
$ws.Cells.Item(24, 11).Value = 70.09999999999999  # K24: was 76.22221999999999
$ws.Cells.Item(24, 13).Value = 99.90000000000001  # M24: was 93.77778000000001
$ws.Cells.Item(107, 8).Value = 598.1429000000001  # H107: was 667
$ws.Cells.Item(107, 9).Value = 214.5  # I107: was 220.4
$ws.Cells.Item(107, 11).Value = 214.5  # K107: was 220.4
$ws.Cells.Item(107, 13).Value = 1705.5  # M107: was 1699.6
$ws.Cells.Item(132, 8).Value = 7248.8184  # H132: was 7965.6665
$ws.Cells.Item(132, 9).Value = 4390.0835  # I132: was 5090.4165
$ws.Cells.Item(132, 10).Value = 10679.3  # J132: was 11799.333
$ws.Cells.Item(132, 11).Value = 13170.2505  # K132: was 15271.2495
$ws.Cells.Item(132, 12).Value = 32037.9  # L132: was 35397.999
$ws.Cells.Item(132, 13).Value = -10640.2505  # M132: was -12741.2495
$ws.Cells.Item(132, 14).Value = -37097.89999999999  # N132: was -40457.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 302  # H5: was 301
$ws.Cells.Item(5, 9).Value = 302.66666  # I5: was 301
$ws.Cells.Item(5, 10).Value = 300  # J5: was 0
$ws.Cells.Item(5, 11).Value = 907.9999799999999  # K5: was 903
$ws.Cells.Item(5, 12).Value = 900  # L5: was 0
$ws.Cells.Item(5, 13).Value = -795.9999799999999  # M5: was -791
$ws.Cells.Item(5, 14).Value = -1124  # N5: was None
$ws.Cells.Item(7, 8).Value = 0  # H7: was 1000
$ws.Cells.Item(7, 10).Value = 0  # J7: was 1000
$ws.Cells.Item(7, 12).Value = 0  # L7: was 3000
$ws.Cells.Item(7, 14).Value = ""  # N7: was -3224
$ws.Cells.Item(8, 8).Value = 326.86667  # H8: was 356.83334
$ws.Cells.Item(8, 9).Value = 326.86667  # I8: was 356.83334
$ws.Cells.Item(8, 11).Value = 980.60001  # K8: was 1070.50002
$ws.Cells.Item(8, 13).Value = -841.60001  # M8: was -931.5000199999999
$ws.Cells.Item(44, 8).Value = 658.1818  # H44: was 691.1429000000001
$ws.Cells.Item(44, 9).Value = 790  # I44: was 691.1429000000001
$ws.Cells.Item(44, 10).Value = 500  # J44: was 0
$ws.Cells.Item(44, 11).Value = 2370  # K44: was 2073.4287
$ws.Cells.Item(44, 12).Value = 1500  # L44: was 0
$ws.Cells.Item(44, 13).Value = -1972  # M44: was -1675.4287
$ws.Cells.Item(44, 14).Value = -2296  # N44: was None
$ws.Cells.Item(47, 8).Value = 158.33333  # H47: was 293.75
$ws.Cells.Item(47, 9).Value = 158.33333  # I47: was 293.75
$ws.Cells.Item(47, 11).Value = 474.99999  # K47: was 881.25
$ws.Cells.Item(47, 13).Value = -43.99998999999997  # M47: was -450.25
$ws.Cells.Item(69, 8).Value = 12  # H69: was 0
$ws.Cells.Item(69, 9).Value = 12  # I69: was 0
$ws.Cells.Item(69, 11).Value = 36  # K69: was 0
$ws.Cells.Item(69, 13).Value = 775  # M69: was None
$ws.Cells.Item(72, 8).Value = 12  # H72: was 0
$ws.Cells.Item(72, 9).Value = 12  # I72: was 0
$ws.Cells.Item(72, 11).Value = 108  # K72: was 0
$ws.Cells.Item(72, 13).Value = 3948  # M72: was None
$ws.Cells.Item(113, 8).Value = 496.5  # H113: was 549.6667
$ws.Cells.Item(113, 9).Value = 0  # I113: was 650
$ws.Cells.Item(113, 10).Value = 496.5  # J113: was 499.5
$ws.Cells.Item(113, 11).Value = 0  # K113: was 1950
$ws.Cells.Item(113, 12).Value = 1489.5  # L113: was 1498.5
$ws.Cells.Item(113, 13).Value = ""  # M113: was 220
$ws.Cells.Item(113, 14).Value = -5829.5  # N113: was -5838.5
$ws.Cells.Item(120, 8).Value = 400  # H120: was 0
$ws.Cells.Item(120, 9).Value = 400  # I120: was 0
$ws.Cells.Item(120, 11).Value = 1200  # K120: was 0
$ws.Cells.Item(120, 13).Value = 3638  # M120: was None
$ws.Cells.Item(128, 8).Value = 499992  # H128: was 0
$ws.Cells.Item(128, 9).Value = 499992  # I128: was 0
$ws.Cells.Item(128, 11).Value = 1499976  # K128: was 0
$ws.Cells.Item(128, 13).Value = -1494996  # M128: was None
$ws.Cells.Item(135, 8).Value = 302  # H135: was 301
$ws.Cells.Item(135, 9).Value = 302.66666  # I135: was 301
$ws.Cells.Item(135, 10).Value = 300  # J135: was 0
$ws.Cells.Item(135, 11).Value = 2723.99994  # K135: was 2709
$ws.Cells.Item(135, 12).Value = 2700  # L135: was 0
$ws.Cells.Item(135, 13).Value = -188.9999399999997  # M135: was -174
$ws.Cells.Item(135, 14).Value = -7770  # N135: was None

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(9, 8).Value = 283.5  # H9: was 333
$ws.Cells.Item(9, 9).Value = 283.5  # I9: was 333
$ws.Cells.Item(9, 11).Value = 283.5  # K9: was 333
$ws.Cells.Item(9, 13).Value = -113.5  # M9: was -163
$ws.Cells.Item(80, 8).Value = 3300  # H80: was 2431.6667
$ws.Cells.Item(80, 9).Value = 3300  # I80: was 2647.5
$ws.Cells.Item(80, 10).Value = 0  # J80: was 2000
$ws.Cells.Item(80, 11).Value = 3300  # K80: was 2647.5
$ws.Cells.Item(80, 12).Value = 0  # L80: was 2000
$ws.Cells.Item(80, 13).Value = -2302  # M80: was -1649.5
$ws.Cells.Item(80, 14).Value = ""  # N80: was -3996
$ws.Cells.Item(83, 8).Value = 3300  # H83: was 2431.6667
$ws.Cells.Item(83, 9).Value = 3300  # I83: was 2647.5
$ws.Cells.Item(83, 10).Value = 0  # J83: was 2000
$ws.Cells.Item(83, 11).Value = 16500  # K83: was 13237.5
$ws.Cells.Item(83, 12).Value = 0  # L83: was 10000
$ws.Cells.Item(83, 13).Value = -11508  # M83: was -8245.5
$ws.Cells.Item(83, 14).Value = ""  # N83: was -19984
$ws.Cells.Item(107, 8).Value = 802.8333  # H107: was 795
$ws.Cells.Item(107, 9).Value = 762  # I107: was 755.9
$ws.Cells.Item(107, 11).Value = 762  # K107: was 755.9
$ws.Cells.Item(107, 13).Value = 1158  # M107: was 1164.1
$ws.Cells.Item(136, 8).Value = 29484.133  # H136: was 29596
$ws.Cells.Item(136, 10).Value = 29484.133  # J136: was 29596
$ws.Cells.Item(136, 12).Value = 88452.399  # L136: was 88788
$ws.Cells.Item(136, 14).Value = -93552.399  # N136: was -93888

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(9, 8).Value = 864.75  # H9: was 371.16666
$ws.Cells.Item(9, 9).Value = 343.8  # I9: was 365.6
$ws.Cells.Item(9, 10).Value = 1733  # J9: was 399
$ws.Cells.Item(9, 11).Value = 343.8  # K9: was 365.6
$ws.Cells.Item(9, 12).Value = 1733  # L9: was 399
$ws.Cells.Item(9, 13).Value = -119.8  # M9: was -141.6
$ws.Cells.Item(9, 14).Value = -2181  # N9: was -847
$ws.Cells.Item(30, 8).Value = 2582.2727  # H30: was 453.85715
$ws.Cells.Item(30, 9).Value = 840.6  # I30: was 453.85715
$ws.Cells.Item(30, 10).Value = 19999  # J30: was 0
$ws.Cells.Item(30, 11).Value = 840.6  # K30: was 453.85715
$ws.Cells.Item(30, 12).Value = 19999  # L30: was 0
$ws.Cells.Item(30, 13).Value = -732.6  # M30: was -345.85715
$ws.Cells.Item(30, 14).Value = -20215  # N30: was None
$ws.Cells.Item(35, 8).Value = 2740.5557  # H35: was 2473.5
$ws.Cells.Item(35, 9).Value = 583.125  # I35: was 526.1111
$ws.Cells.Item(35, 11).Value = 583.125  # K35: was 526.1111
$ws.Cells.Item(35, 13).Value = -247.125  # M35: was -190.1111
$ws.Cells.Item(46, 8).Value = 4613.2666  # H46: was 4406.1875
$ws.Cells.Item(46, 10).Value = 4777.6665  # J46: was 4429.9
$ws.Cells.Item(46, 12).Value = 4777.6665  # L46: was 4429.9
$ws.Cells.Item(46, 14).Value = -5153.6665  # N46: was -4805.9
$ws.Cells.Item(61, 8).Value = 1400  # H61: was 1700
$ws.Cells.Item(61, 10).Value = 0  # J61: was 2000
$ws.Cells.Item(61, 12).Value = 0  # L61: was 2000
$ws.Cells.Item(61, 14).Value = ""  # N61: was -2404
$ws.Cells.Item(62, 8).Value = 0  # H62: was 32499.5
$ws.Cells.Item(62, 9).Value = 0  # I62: was 5000
$ws.Cells.Item(62, 10).Value = 0  # J62: was 59999
$ws.Cells.Item(62, 11).Value = 0  # K62: was 5000
$ws.Cells.Item(62, 12).Value = 0  # L62: was 59999
$ws.Cells.Item(62, 13).Value = ""  # M62: was -4376
$ws.Cells.Item(62, 14).Value = ""  # N62: was -61247
$ws.Cells.Item(65, 8).Value = 0  # H65: was 32499.5
$ws.Cells.Item(65, 9).Value = 0  # I65: was 5000
$ws.Cells.Item(65, 10).Value = 0  # J65: was 59999
$ws.Cells.Item(65, 11).Value = 0  # K65: was 15000
$ws.Cells.Item(65, 12).Value = 0  # L65: was 179997
$ws.Cells.Item(65, 13).Value = ""  # M65: was -11880
$ws.Cells.Item(65, 14).Value = ""  # N65: was -186237
$ws.Cells.Item(96, 8).Value = 40000  # H96: was 36666.668
$ws.Cells.Item(96, 10).Value = 40000  # J96: was 36666.668
$ws.Cells.Item(96, 12).Value = 40000  # L96: was 36666.668
$ws.Cells.Item(96, 14).Value = -45492  # N96: was -42158.668
$ws.Cells.Item(100, 8).Value = 2000  # H100: was 1295
$ws.Cells.Item(100, 9).Value = 2000  # I100: was 1295
$ws.Cells.Item(100, 11).Value = 2000  # K100: was 1295
$ws.Cells.Item(100, 13).Value = -1459  # M100: was -754
$ws.Cells.Item(106, 8).Value = 39857.145  # H106: was 0
$ws.Cells.Item(106, 10).Value = 39857.145  # J106: was 0
$ws.Cells.Item(106, 12).Value = 39857.145  # L106: was 0
$ws.Cells.Item(106, 14).Value = -42381.145  # N106: was None
$ws.Cells.Item(113, 8).Value = 1400  # H113: was 1700
$ws.Cells.Item(113, 10).Value = 0  # J113: was 2000
$ws.Cells.Item(113, 12).Value = 0  # L113: was 2000
$ws.Cells.Item(113, 14).Value = ""  # N113: was -6340
$ws.Cells.Item(136, 8).Value = 4447.5  # H136: was 4515.875
$ws.Cells.Item(136, 9).Value = 3808.25  # I136: was 3859.4546
$ws.Cells.Item(136, 10).Value = 5299.8335  # J136: was 5960
$ws.Cells.Item(136, 11).Value = 11424.75  # K136: was 11578.3638
$ws.Cells.Item(136, 12).Value = 15899.5005  # L136: was 17880
$ws.Cells.Item(136, 13).Value = -8874.75  # M136: was -9028.363799999999
$ws.Cells.Item(136, 14).Value = -20999.5005  # N136: was -22980

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(17, 8).Value = 683.3333  # H17: was 180
$ws.Cells.Item(17, 9).Value = 683.3333  # I17: was 180
$ws.Cells.Item(17, 11).Value = 683.3333  # K17: was 180
$ws.Cells.Item(17, 13).Value = -511.3333  # M17: was -8
$ws.Cells.Item(34, 8).Value = 0  # H34: was 8513
$ws.Cells.Item(34, 9).Value = 0  # I34: was 2026
$ws.Cells.Item(34, 10).Value = 0  # J34: was 15000
$ws.Cells.Item(34, 11).Value = 0  # K34: was 2026
$ws.Cells.Item(34, 12).Value = 0  # L34: was 15000
$ws.Cells.Item(34, 13).Value = ""  # M34: was -1823
$ws.Cells.Item(34, 14).Value = ""  # N34: was -15406
$ws.Cells.Item(132, 8).Value = 2919.2  # H132: was 2781.8333
$ws.Cells.Item(132, 9).Value = 3149  # I132: was 3174
$ws.Cells.Item(132, 10).Value = 2000  # J132: was 1997.5
$ws.Cells.Item(132, 11).Value = 9447  # K132: was 9522
$ws.Cells.Item(132, 12).Value = 6000  # L132: was 5992.5
$ws.Cells.Item(132, 13).Value = -6917  # M132: was -6992
$ws.Cells.Item(132, 14).Value = -11060  # N132: was -11052.5
